$wb = $excel.ActiveWorkbook

# --- Sheet1 (monsoon19): update the selection, no longer the active tab ---
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "passwords" sheet after monsoon19 ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "passwords"

# Write values in this order so the shared-string table gets
# "not set" before "stored" (matches index 15 / 16 ordering).
$newSheet.Range("B2").Value = "not set"
$newSheet.Range("B3").Value = "not set"
$newSheet.Range("B1").Value = "stored"

$newSheet.Range("A2").Value = "archit.checker_ug20"
$newSheet.Range("A3").Value = "rathi.kashi_ug20"
$newSheet.Range("A4").Value = "deepraj.pandey_ug20"
$newSheet.Range("A5").Value = "aastha.shah_ug20"
$newSheet.Range("A6").Value = "reuel.john_ug20"
$newSheet.Range("A7").Value = "yash.dixit_ug20"

$newSheet.Columns.Item(1).ColumnWidth = 18

# Set the selection/active cell on the new sheet and make it the active tab
$newSheet.Range("B1").Select()

# Now fix up the selection on sheet1 (it is no longer the active tab)
$ws1.Range("A2:A7").Select()
$newSheet.Activate()
